# Method of Hierarchy Analysis — add the "Дизайн" criterion label to
# the matrix's corner cell and horizontally center the pairwise
# comparison values / priority vector.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corner cell (A1): label the criterion. It keeps the Times New Roman
# font + thin border already shared by the header/label cells, but
# (being a short single word) doesn't need the wrap/vertical-centering
# alignment those cells use, so that gets cleared back to the default.
$ws.Range("A1").Value = "Дизайн"
$ws.Range("A1").WrapText = $false
$ws.Range("A1").VerticalAlignment = -4107   # xlBottom (default)

# Pairwise-comparison matrix + priority-vector column (B2:H7):
# horizontally center the values in addition to the vertical
# centering / wrap they already have.
$ws.Range("B2:H7").HorizontalAlignment = -4108   # xlCenter

Write-Host "Applied corner label and centered the comparison matrix."
